$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "26.659.37"
Set-TextValue "E2" "  +0.78%  "
Set-TextValue "D3" "1.644.58"
Set-TextValue "E3" "  +1.09%  "
Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.17%  "
Set-TextValue "D5" "215.44"
Set-TextValue "E5" "  +1.10%  "
Set-TextValue "D6" "0.506"
Set-TextValue "E6" "  +1.32%  "
Set-TextValue "E7" "  +0.13%  "
Set-TextValue "D8" "0.251"
Set-TextValue "E8" "  +1.14%  "
Set-TextValue "D9" "0.0625"
Set-TextValue "E9" "  +0.20%  "
Set-TextValue "D10" "19.12"
Set-TextValue "E10" "  +1.22%  "
Set-TextValue "D11" "0.0843"
Set-TextValue "E11" "  -0.20%  "
Set-TextValue "D12" "1.877.25"
Set-TextValue "E12" "  +1.26%  "
Set-TextValue "D13" "1.699.92"
Set-TextValue "E13" "  +4.54%  "
Set-TextValue "D14" "4.17"
Set-TextValue "E14" "  +0.98%  "
Set-TextValue "D15" "0.531"
Set-TextValue "E15" "  +1.61%  "
Set-TextValue "D16" "65.05"
Set-TextValue "E16" "  +0.36%  "
Set-TextValue "D17" "26.684.77"
Set-TextValue "E17" "  +0.63%  "
Set-TextValue "D18" "0.0₃0741"
Set-TextValue "E18" "  +0.18%  "
Set-TextValue "D19" "217.93"
Set-TextValue "E19" "  +1.35%  "
Set-TextValue "E20" "  +0.17%  "
Set-TextValue "D21" "4.35"
Set-TextValue "E21" "  +1.11%  "
Set-TextValue "D22" "6.26"
Set-TextValue "E22" "  -0.11%  "
Set-TextValue "D23" "9.49"
Set-TextValue "E23" "  +2.24%  "
Set-TextValue "E24" "  +13.43%  "
Set-TextValue "D25" "145.67"
Set-TextValue "E25" "  -1.92%  "
Set-TextValue "E26" "  +0.18%  "
Set-TextValue "E27" "  +0.35%  "
Set-TextValue "D28" "7.11"
Set-TextValue "E28" "  +4.15%  "
Set-TextValue "D29" "15.70"
Set-TextValue "E29" "  +0.89%  "
Set-TextValue "E30" "  +1.24%  "
Set-TextValue "D31" "1.17"
Set-TextValue "E31" "  +1.03%  "
Set-TextValue "D32" "3.36"
Set-TextValue "E32" "  +0.99%  "
Set-TextValue "D33" "3.00"
Set-TextValue "E33" "  +1.89%  "
Set-TextValue "D34" "1.275.43"
Set-TextValue "E34" "  +4.63%  "
Set-TextValue "E35" "  +3.36%  "
Set-TextValue "D36" "2.42"
Set-TextValue "E36" "  +1.41%  "
Set-TextValue "D37" "0.0178"
Set-TextValue "E37" "  +2.61%  "
Set-TextValue "E38" "  +5.55%  "
Set-TextValue "D39" "0.822"
Set-TextValue "E39" "  +3.58%  "
Set-TextValue "D40" "1.00"
Set-TextValue "E40" "  -0.04%  "
Set-TextValue "D41" "0.812"
Set-TextValue "E41" "  +2.49%  "
Set-TextValue "E42" "  -0.95%  "
Set-TextValue "D43" "5.44"
Set-TextValue "E43" "  +1.40%  "
Set-TextValue "D44" "1.786.76"
Set-TextValue "E44" "  +1.22%  "
Set-TextValue "D45" "91.66"
Set-TextValue "E45" "  -1.57%  "
Set-TextValue "D46" "59.82"
Set-TextValue "E46" "  +9.06%  "
Set-TextValue "E47" "  +1.26%  "
Set-TextValue "E48" "  +1.10%  "
Set-TextValue "D49" "7.78"
Set-TextValue "E49" "  +3.18%  "
Set-TextValue "D50" "0.0967"
Set-TextValue "E50" "  +1.81%  "
Set-TextValue "D51" "0.408"
Set-TextValue "E51" "  +0.29%  "
